$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 41) with a new record, mirroring the pattern of row 40
$ws.Range("A41").Value = "PB"
$ws.Range("B41").Value = "ZRP803_CHD_P45"
$ws.Range("C41").Value = "L1800"
$ws.Range("D41").Value = "18-Dec-2025 6:27 PM"
$ws.Range("E41").Value = "FAIL"
$ws.Range("F41").Value = "1. Network Detach Success rate"
$ws.Range("G41").Value = "1. Static Cell Reselection"
$ws.Range("H41").Value = "1. Exclude the current logfile and create a new one. Ensure that the Detach Request matches the Detach Accept. Verify this in the Event tab of AZQ before saving the logfile"

# Copy formatting from row 40 to row 41 (same style pattern as prior rows)
$ws.Range("A40:C40").Copy()
$ws.Range("A41:C41").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E40:H40").Copy()
$ws.Range("E41:H41").PasteSpecial(-4122)  # xlPasteFormats

# D41 holds plain text (not a date serial), so it takes the same style as A41/C41
$ws.Range("A41").Copy()
$ws.Range("D41").PasteSpecial(-4122)  # xlPasteFormats

# Set row height to match target (19, a single-line row instead of the wrapped 38)
$ws.Rows.Item(41).RowHeight = 19

# Update selection to match diff's final active cell
$ws.Range("C44").Select()
